$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 4500
$ws.Range("I16").Value = 4500
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 4500
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -4270
$ws.Range("N16").ClearContents()
$ws.Range("H62").Value = 2224
$ws.Range("I62").Value = 2433.3333
$ws.Range("J62").Value = 2014.6666
$ws.Range("K62").Value = 2433.3333
$ws.Range("L62").Value = 2014.6666
$ws.Range("M62").Value = -1809.3333
$ws.Range("N62").Value = -3262.6666
$ws.Range("H65").Value = 2224
$ws.Range("I65").Value = 2433.3333
$ws.Range("J65").Value = 2014.6666
$ws.Range("K65").Value = 12166.6665
$ws.Range("L65").Value = 10073.333
$ws.Range("M65").Value = -9046.666499999999
$ws.Range("N65").Value = -16313.333
$ws.Range("H74").Value = 3354.1516
$ws.Range("I74").Value = 3843.9092
$ws.Range("J74").Value = 3109.2727
$ws.Range("K74").Value = 3843.9092
$ws.Range("L74").Value = 3109.2727
$ws.Range("M74").Value = -2907.9092
$ws.Range("N74").Value = -4981.2727
$ws.Range("H77").Value = 3354.1516
$ws.Range("I77").Value = 3843.9092
$ws.Range("J77").Value = 3109.2727
$ws.Range("K77").Value = 19219.546
$ws.Range("L77").Value = 15546.3635
$ws.Range("M77").Value = -14539.546
$ws.Range("N77").Value = -24906.3635
$ws.Range("H120").Value = 27190.25
$ws.Range("J120").Value = 27190.25
$ws.Range("L120").Value = 27190.25
$ws.Range("N120").Value = -36866.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H25").Value = 1265
$ws.Range("I25").Value = 1265
$ws.Range("K25").Value = 1265
$ws.Range("M25").Value = -863
$ws.Range("H35").Value = 2518.5
$ws.Range("I35").Value = 2518.5
$ws.Range("K35").Value = 2518.5
$ws.Range("M35").Value = -2112.5
$ws.Range("H45").Value = 1475.8182
$ws.Range("I45").Value = 1553.7778
$ws.Range("J45").Value = 1125
$ws.Range("K45").Value = 1553.7778
$ws.Range("L45").Value = 1125
$ws.Range("M45").Value = -1176.7778
$ws.Range("N45").Value = -1879

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 2216.8
$ws.Range("I8").Value = 1521
$ws.Range("K8").Value = 1521
$ws.Range("M8").Value = -1381
$ws.Range("H10").Value = 3000
$ws.Range("J10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("N10").Value = -3280
$ws.Range("H12").Value = 461
$ws.Range("I12").Value = 461
$ws.Range("K12").Value = 461
$ws.Range("M12").Value = -293
$ws.Range("H24").Value = 1032.6666
$ws.Range("I24").Value = 779
$ws.Range("J24").Value = 1540
$ws.Range("K24").Value = 779
$ws.Range("L24").Value = 1540
$ws.Range("M24").Value = -544
$ws.Range("N24").Value = -2010
$ws.Range("H25").Value = 1000
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -765
$ws.Range("N25").ClearContents()
$ws.Range("H29").Value = 5430.3335
$ws.Range("I29").Value = 916.4
$ws.Range("J29").Value = 28000
$ws.Range("K29").Value = 916.4
$ws.Range("L29").Value = 28000
$ws.Range("M29").Value = -627.4
$ws.Range("N29").Value = -28578
$ws.Range("H36").Value = 1504
$ws.Range("I36").Value = 1448.8889
$ws.Range("J36").Value = 2000
$ws.Range("K36").Value = 1448.8889
$ws.Range("L36").Value = 2000
$ws.Range("M36").Value = -914.8888999999999
$ws.Range("N36").Value = -3068
$ws.Range("H37").Value = 826.6667
$ws.Range("I37").Value = 430
$ws.Range("J37").Value = 1025
$ws.Range("K37").Value = 430
$ws.Range("L37").Value = 1025
$ws.Range("M37").Value = -293
$ws.Range("N37").Value = -1299
$ws.Range("H54").Value = 10517.667
$ws.Range("I54").Value = 1237
$ws.Range("J54").Value = 43000
$ws.Range("K54").Value = 1237
$ws.Range("L54").Value = 43000
$ws.Range("M54").Value = -753
$ws.Range("N54").Value = -43968
$ws.Range("H62").Value = 31190
$ws.Range("J62").Value = 31190
$ws.Range("L62").Value = 31190
$ws.Range("N62").Value = -32562
$ws.Range("H65").Value = 31190
$ws.Range("J65").Value = 31190
$ws.Range("L65").Value = 93570
$ws.Range("N65").Value = -100434
$ws.Range("H75").Value = 7347
$ws.Range("I75").Value = 7347
$ws.Range("K75").Value = 7347
$ws.Range("M75").Value = -6411
$ws.Range("H78").Value = 7347
$ws.Range("I78").Value = 7347
$ws.Range("K78").Value = 22041
$ws.Range("M78").Value = -17361

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 2126
$ws.Range("I12").Value = 2126
$ws.Range("K12").Value = 2126
$ws.Range("M12").Value = -1956

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 520.4
$ws.Range("J17").Value = 734
$ws.Range("L17").Value = 2202
$ws.Range("N17").Value = -2540
$ws.Range("H34").Value = 538.9143
$ws.Range("I34").Value = 129
$ws.Range("J34").Value = 972.94116
$ws.Range("K34").Value = 387
$ws.Range("L34").Value = 2918.82348
$ws.Range("M34").Value = -303
$ws.Range("N34").Value = -3086.82348
$ws.Range("H39").Value = 1839.973
$ws.Range("J39").Value = 1914.0883
$ws.Range("L39").Value = 5742.2649
$ws.Range("N39").Value = -6330.2649
$ws.Range("H55").Value = 1207.0588
$ws.Range("I55").Value = 283.33334
$ws.Range("J55").Value = 1405
$ws.Range("K55").Value = 850.0000200000001
$ws.Range("L55").Value = 4215
$ws.Range("M55").Value = -673.0000200000001
$ws.Range("N55").Value = -4569
$ws.Range("H131").Value = 1283034.4
$ws.Range("J131").Value = 1273.8857
$ws.Range("L131").Value = 3821.6571
$ws.Range("N131").Value = -13901.6571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 6250
$ws.Range("J5").Value = 6250
$ws.Range("L5").Value = 6250
$ws.Range("N5").Value = -6474
$ws.Range("H13").Value = 538.125
$ws.Range("I13").Value = 199.5
$ws.Range("J13").Value = 651
$ws.Range("K13").Value = 199.5
$ws.Range("L13").Value = 651
$ws.Range("M13").Value = -60.5
$ws.Range("N13").Value = -929
$ws.Range("H41").Value = 3762.5
$ws.Range("I41").Value = 2525
$ws.Range("J41").Value = 5000
$ws.Range("K41").Value = 2525
$ws.Range("L41").Value = 5000
$ws.Range("M41").Value = -2170
$ws.Range("N41").Value = -5710

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 945
$ws.Range("I16").Value = 945
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 945
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -775
$ws.Range("N16").ClearContents()
$ws.Range("H46").Value = 917.875
$ws.Range("J46").Value = 932.4
$ws.Range("L46").Value = 932.4
$ws.Range("N46").Value = -1308.4
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("N56").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1001500
$ws.Range("J5").Value = 3000
$ws.Range("L5").Value = 3000
$ws.Range("N5").Value = -3224
